# Insert a new data row before the current row 15 (everything from row 15
# down, including the old row 15, shifts down by one row). This brings the
# sheet from A1:T56 to A1:T57, matching the new dimension in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new record's data.
$ws.Range("A15").Value = 5
$ws.Range("B15").Value = "Macroferia Regional de Talca"
$ws.Range("C15").Value = "Maule"
$ws.Range("D15").Value = 45260
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100103
$ws.Range("H15").Value = "Frutos de hueso (carozo)"
$ws.Range("I15").Value = 100103003
$ws.Range("J15").Value = "Damasco"
$ws.Range("K15").Value = "Dina"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 23000
$ws.Range("O15").Value = 23000
$ws.Range("P15").Value = 23000
$ws.Range("Q15").Value = "$/caja 10 kilos"
$ws.Range("R15").Value = "Provincia de Limarí"
$ws.Range("S15").Value = 2300
$ws.Range("T15").Value = 10
